$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep a reference to a plain (unstyled) cell style so that writing numeric-looking
# text into column D does not pick up an implicit Text number format / style index,
# matching the original inlineStr cells which carry no style attribute.
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.941.93"
$ws.Range("D2").Style = $plainStyle
$ws.Range("E2").Value = "  -0.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.523.76"
$ws.Range("D3").Style = $plainStyle
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.78"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  +0.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.11"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  -0.89%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = "  +0.81%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.522.31"
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = "  -0.22%  "

$ws.Range("E10").Value = "  +0.88%  "

$ws.Range("E11").Value = "  -1.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.36"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  -1.42%  "

$ws.Range("E13").Value = "  -1.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.971.91"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = "  +0.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.16"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = "  -0.70%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.956.70"
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = "  -0.31%  "

$ws.Range("E17").Value = "  -0.84%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.523.85"
$ws.Range("D18").Style = $plainStyle
$ws.Range("E18").Value = "  +0.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.12"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = "  +0.70%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.29"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = "  -0.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.68"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  +0.37%  "

$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.98"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "  +2.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.97"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  +5.73%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.425"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = "  -0.55%  "

$ws.Range("E26").Value = "  +1.21%  "

$ws.Range("E27").Value = "  +0.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.67"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = "  -1.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.74"
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = "  -0.46%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0776"
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = "  +0.39%  "

$ws.Range("E31").Value = "  +0.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.19"
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = "  +5.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "162.24"
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = "  -1.66%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.49"
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = "  +1.17%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = $plainStyle

$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.14"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = "  -2.96%  "

$ws.Range("E38").Value = "  -1.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.64"
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = "  -0.71%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.819"
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = "  +1.25%  "

$ws.Range("E41").Value = "  -0.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "287.53"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = "  +2.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.22"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = "  -0.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "132.42"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  +8.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.997"
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.611"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "  +2.18%  "

$ws.Range("E47").Value = "  +0.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0934"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  +0.09%  "

$ws.Range("E49").Value = "  -0.34%  "

$ws.Range("E50").Value = "  -1.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.42"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "  -1.86%  "

